# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# For rows 16..53 of "Hoja1" (the "Periodo Mora" table):
#  - Column E holds the period label (e.g. "2003", "1701", ...). The table
#    is re-sorted so the periods run in ascending order (1701 -> 2003)
#    instead of the previous descending order.
#  - Column F ("Valor Mora") and column G ("Salario Basico") are refreshed
#    with the new values coming from the updated base de datos.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = [ordered]@{
    16 = @{ Period = "1701"; F = 27578; G = 781242 }
    17 = @{ Period = "1702"; F = 27578; G = 781242 }
    18 = @{ Period = "1704"; F = 27578; G = 781242 }
    19 = @{ Period = "1705"; F = 27578; G = 781242 }
    20 = @{ Period = "1706"; F = 27578; G = 781242 }
    21 = @{ Period = "1707"; F = 27578; G = 781242 }
    22 = @{ Period = "1708"; F = 27578; G = 781242 }
    23 = @{ Period = "1709"; F = 27578; G = 781242 }
    24 = @{ Period = "1710"; F = 27578; G = 781242 }
    25 = @{ Period = "1711"; F = 27578; G = 781242 }
    26 = @{ Period = "1712"; F = 27578; G = 781242 }
    27 = @{ Period = "1801"; F = 27578; G = 781242 }
    28 = @{ Period = "1802"; F = 27578; G = 781242 }
    29 = @{ Period = "1803"; F = 27578; G = 781242 }
    30 = @{ Period = "1804"; F = 27578; G = 781242 }
    31 = @{ Period = "1805"; F = 27578; G = 781242 }
    32 = @{ Period = "1806"; F = 27578; G = 781242 }
    33 = @{ Period = "1807"; F = 27578; G = 781242 }
    34 = @{ Period = "1808"; F = 27578; G = 781242 }
    35 = @{ Period = "1809"; F = 31249; G = 781242 }
    36 = @{ Period = "1810"; F = 31249; G = 781242 }
    37 = @{ Period = "1811"; F = 31249; G = 781242 }
    38 = @{ Period = "1812"; F = 31249; G = 781242 }
    39 = @{ Period = "1901"; F = 31249; G = 781242 }
    40 = @{ Period = "1902"; F = 31249; G = 781242 }
    41 = @{ Period = "1903"; F = 31249; G = 781242 }
    42 = @{ Period = "1904"; F = 31249; G = 781242 }
    43 = @{ Period = "1905"; F = 31249; G = 781242 }
    44 = @{ Period = "1906"; F = 31249; G = 781242 }
    45 = @{ Period = "1907"; F = 31249; G = 781242 }
    46 = @{ Period = "1908"; F = 31249; G = 781242 }
    47 = @{ Period = "1909"; F = 31249; G = 781242 }
    48 = @{ Period = "1910"; F = 31249; G = 781242 }
    49 = @{ Period = "1911"; F = 31249; G = 781242 }
    50 = @{ Period = "1912"; F = 31249; G = 781242 }
    51 = @{ Period = "2001"; F = 31249; G = 781242 }
    52 = @{ Period = "2002"; F = 31249; G = 781242 }
    53 = @{ Period = "2003"; F = 30208; G = 781242 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Cells.Item($r, 5).Value = $data.Period
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
}
